$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 12896.556
$ws.Range("J17").Value = 12896.556
$ws.Range("L17").Value = 38689.66800000001
$ws.Range("N17").Value = -39025.66800000001

$ws.Range("H28").Value = 656.7
$ws.Range("I28").Value = 388.8125
$ws.Range("J28").Value = 1728.25
$ws.Range("K28").Value = 388.8125
$ws.Range("L28").Value = 1728.25
$ws.Range("M28").Value = 96.1875
$ws.Range("N28").Value = -2698.25

$ws.Range("H113").Value = 5116.1665
$ws.Range("J113").Value = 6481.4
$ws.Range("L113").Value = 6481.4
$ws.Range("N113").Value = -12989.4

$ws.Range("H132").Value = 411.5
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""

$ws.Range("H135").Value = 864.3214
$ws.Range("I135").Value = 567.5454999999999
$ws.Range("J135").Value = 1952.5
$ws.Range("K135").Value = 5107.9095
$ws.Range("L135").Value = 17572.5
$ws.Range("M135").Value = -2572.9095
$ws.Range("N135").Value = -22642.5

$ws.Range("H137").Value = 35935.02
$ws.Range("I137").Value = 43758.094
$ws.Range("J137").Value = 3078.1
$ws.Range("K137").Value = 131274.282
$ws.Range("L137").Value = 9234.299999999999
$ws.Range("M137").Value = -128724.282
$ws.Range("N137").Value = -14334.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13506.909
$ws.Range("I32").Value = 8526.344999999999
$ws.Range("J32").Value = 23136
$ws.Range("K32").Value = 8526.344999999999
$ws.Range("L32").Value = 23136
$ws.Range("M32").Value = -8239.344999999999
$ws.Range("N32").Value = -23710

$ws.Range("H45").Value = 8101461
$ws.Range("I45").Value = 12823231
$ws.Range("J45").Value = 6998.2856
$ws.Range("K45").Value = 12823231
$ws.Range("L45").Value = 6998.2856
$ws.Range("M45").Value = -12822854
$ws.Range("N45").Value = -7752.2856

$ws.Range("H61").Value = 5129.4707
$ws.Range("I61").Value = 5147.4
$ws.Range("K61").Value = 5147.4
$ws.Range("M61").Value = -4935.4

$ws.Range("H74").Value = 23683.38
$ws.Range("J74").Value = 59610.25
$ws.Range("L74").Value = 59610.25
$ws.Range("N74").Value = -61358.25

$ws.Range("H77").Value = 23683.38
$ws.Range("J77").Value = 59610.25
$ws.Range("L77").Value = 298051.25
$ws.Range("N77").Value = -306787.25

$ws.Range("H122").Value = 1162720
$ws.Range("I122").Value = 3523.4285
$ws.Range("J122").Value = 1900390.6
$ws.Range("K122").Value = 10570.2855
$ws.Range("L122").Value = 5701171.800000001
$ws.Range("M122").Value = -8120.2855
$ws.Range("N122").Value = -5706071.800000001

$ws.Range("H132").Value = 2314.3618
$ws.Range("I132").Value = 1784.8605
$ws.Range("K132").Value = 5354.5815
$ws.Range("M132").Value = -2824.5815

$ws.Range("H136").Value = 5129.4707
$ws.Range("I136").Value = 5147.4
$ws.Range("K136").Value = 15442.2
$ws.Range("M136").Value = -12892.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 16669479
$ws.Range("I20").Value = 41669264
$ws.Range("J20").Value = 2954.8333
$ws.Range("K20").Value = 41669264
$ws.Range("L20").Value = 2954.8333
$ws.Range("M20").Value = -41669017
$ws.Range("N20").Value = -3448.8333

$ws.Range("H105").Value = 12500900
$ws.Range("I105").Value = 15625825
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 15625825
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = -15624078
$ws.Range("N105").Value = -4694

$ws.Range("H107").Value = 2859701.2
$ws.Range("I107").Value = 3403596.8
$ws.Range("K107").Value = 3403596.8
$ws.Range("M107").Value = -3401676.8

$ws.Range("H123").Value = 24999
$ws.Range("J123").Value = 24999
$ws.Range("L123").Value = 24999
$ws.Range("N123").Value = -34799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24511.6
$ws.Range("I31").Value = 2625.5386
$ws.Range("J31").Value = 33402.812
$ws.Range("K31").Value = 2625.5386
$ws.Range("L31").Value = 33402.812
$ws.Range("M31").Value = -2330.5386
$ws.Range("N31").Value = -33992.812

$ws.Range("H34").Value = 24511.6
$ws.Range("I34").Value = 2625.5386
$ws.Range("J34").Value = 33402.812
$ws.Range("K34").Value = 2625.5386
$ws.Range("L34").Value = 33402.812
$ws.Range("M34").Value = -2423.5386
$ws.Range("N34").Value = -33806.812

$ws.Range("H122").Value = 2754.7896
$ws.Range("I122").Value = 2630.0557
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7890.1671
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -5440.1671
$ws.Range("N122").Value = -19900

$ws.Range("H134").Value = 3851.16
$ws.Range("I134").Value = 2546.0588
$ws.Range("K134").Value = 7638.176399999999
$ws.Range("M134").Value = -5103.176399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 453.33334
$ws.Range("J86").Value = 1499
$ws.Range("L86").Value = 4497
$ws.Range("N86").Value = -6869

$ws.Range("H89").Value = 453.33334
$ws.Range("J89").Value = 1499
$ws.Range("L89").Value = 13491
$ws.Range("N89").Value = -25347

$ws.Range("H131").Value = 23153676
$ws.Range("J131").Value = 22227610
$ws.Range("L131").Value = 66682830
$ws.Range("N131").Value = -66692910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 40004900
$ws.Range("I70").Value = 50004376
$ws.Range("K70").Value = 50004376
$ws.Range("M70").Value = -50004106

$ws.Range("H73").Value = 40004900
$ws.Range("I73").Value = 50004376
$ws.Range("K73").Value = 50004376
$ws.Range("M73").Value = -50003440

$ws.Range("H102").Value = 6179961.5
$ws.Range("I102").Value = 7408581.5
$ws.Range("J102").Value = 3547203.2
$ws.Range("K102").Value = 7408581.5
$ws.Range("L102").Value = 3547203.2
$ws.Range("M102").Value = -7406959.5
$ws.Range("N102").Value = -3550447.2

$ws.Range("H126").Value = 5857868.5
$ws.Range("I126").Value = 5684499
$ws.Range("J126").Value = 5956937
$ws.Range("K126").Value = 17053497
$ws.Range("L126").Value = 17870811
$ws.Range("M126").Value = -17051027
$ws.Range("N126").Value = -17875751

$ws.Range("H132").Value = 2908.3823
$ws.Range("I132").Value = 2673.8708
$ws.Range("J132").Value = 5331.6665
$ws.Range("K132").Value = 8021.6124
$ws.Range("L132").Value = 15994.9995
$ws.Range("M132").Value = -5491.6124
$ws.Range("N132").Value = -21054.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7070.1
$ws.Range("I40").Value = 4284
$ws.Range("K40").Value = 4284
$ws.Range("M40").Value = -4148

$ws.Range("H82").Value = 211114540
$ws.Range("I82").Value = 211114540
$ws.Range("K82").Value = 211114540
$ws.Range("M82").Value = -211114179

$ws.Range("H85").Value = 211114540
$ws.Range("I85").Value = 211114540
$ws.Range("K85").Value = 211114540
$ws.Range("M85").Value = -211113292

$ws.Range("H136").Value = 42010.85
$ws.Range("I136").Value = 65705.875
$ws.Range("J136").Value = 7545.364
$ws.Range("K136").Value = 197117.625
$ws.Range("L136").Value = 22636.092
$ws.Range("M136").Value = -194567.625
$ws.Range("N136").Value = -27736.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13894647
$ws.Range("I81").Value = 41669492
$ws.Range("J81").Value = 7223.5
$ws.Range("K81").Value = 83338984
$ws.Range("L81").Value = 14447
$ws.Range("M81").Value = -83337923
$ws.Range("N81").Value = -16569

$ws.Range("H84").Value = 13894647
$ws.Range("I84").Value = 41669492
$ws.Range("J84").Value = 7223.5
$ws.Range("K84").Value = 416694920
$ws.Range("L84").Value = 72235
$ws.Range("M84").Value = -416689616
$ws.Range("N84").Value = -82843

$ws.Range("H136").Value = 3637.027
$ws.Range("I136").Value = 3314.7585
$ws.Range("J136").Value = 4805.25
$ws.Range("K136").Value = 9944.2755
$ws.Range("L136").Value = 14415.75
$ws.Range("M136").Value = -7394.2755
$ws.Range("N136").Value = -19515.75
